# Apply the commit's change to the tasks workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new comment line to the existing Issues/Comments text in D16.
$existing = $ws.Range("D16").Value2
$ws.Range("D16").Value2 = $existing + "`nLight gbm gave an issue so I had to install some other lightgb library as root and added to the docker file"

# Update the row height for row 16 to fit the new text (90 -> 135).
$ws.Rows.Item(16).RowHeight = 135

# Update the selection so the active cell in the frozen pane reflects D16.
$ws.Activate()
$ws.Range("D16").Select()
